$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.219.25"
$ws.Range("E2").Value = "  +3.17%  "

$ws.Range("D3").Value = "2.574.51"
$ws.Range("E3").Value = "  +4.52%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'505.34"
$ws.Range("E5").Value = "  +2.12%  "

$ws.Range("D6").Value = "'153.81"
$ws.Range("E6").Value = "  -3.42%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  -5.21%  "

$ws.Range("D9").Value = "2.586.24"
$ws.Range("E9").Value = "  +3.65%  "

$ws.Range("D10").Value = "'6.59"
$ws.Range("E10").Value = "  +4.04%  "

$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("D12").Value = "'0.342"
$ws.Range("E12").Value = "  +1.53%  "

$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").Value = "3.019.26"
$ws.Range("E14").Value = "  +4.45%  "

$ws.Range("D15").Value = "60.203.24"
$ws.Range("E15").Value = "  +3.49%  "

$ws.Range("D16").Value = "'21.60"
$ws.Range("E16").Value = "  +1.32%  "

$ws.Range("E17").Value = "  +2.87%  "

$ws.Range("D18").Value = "2.578.14"
$ws.Range("E18").Value = "  +4.00%  "

$ws.Range("D19").Value = "'4.80"
$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("D20").Value = "'345.33"
$ws.Range("E20").Value = "  +5.05%  "

$ws.Range("D21").Value = "'10.32"
$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").Value = "'6.05"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").Value = "'60.07"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").Value = "'0.420"
$ws.Range("E25").Value = "  +2.67%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.164"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.679.14"
$ws.Range("E27").Value = "  +4.69%  "

$ws.Range("D28").Value = "'0.992"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("D29").Value = "0.0₃0850"
$ws.Range("E29").Value = "  +5.20%  "

$ws.Range("D30").Value = "'7.45"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").Value = "'155.72"
$ws.Range("E32").Value = "  +2.91%  "

$ws.Range("D33").Value = "'19.17"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("D35").Value = "'5.73"
$ws.Range("E35").Value = "  +5.67%  "

$ws.Range("D36").Value = "'3.98"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("E37").Value = "  +2.31%  "

$ws.Range("D38").Value = "'0.859"
$ws.Range("E38").Value = "  +22.16%  "

$ws.Range("D39").Value = "'0.846"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  +3.60%  "

$ws.Range("D41").Value = "'3.75"
$ws.Range("E41").Value = "  +3.82%  "

$ws.Range("D42").Value = "'300.11"
$ws.Range("E42").Value = "  +6.04%  "

$ws.Range("D43").Value = "'35.55"
$ws.Range("E43").Value = "  +3.18%  "

$ws.Range("D44").Value = "'0.0565"
$ws.Range("E44").Value = "  +3.71%  "

$ws.Range("D45").Value = "'0.0997"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").Value = "'0.615"
$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("D48").Value = "'19.76"
$ws.Range("E48").Value = "  +8.43%  "

$ws.Range("D49").Value = "'4.95"
$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0233"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.026.29"
$ws.Range("E51").Value = "  +6.31%  "
